$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (the merged "Area of ROI (1280*1024)" title row),
# which shifts everything up by one row.
$ws.Rows("1").Delete()

# Scroll / selection state to match the recorded view.
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("K1:K1048576").Select()
